$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers
$ws.Range("F1").Value = "status"
$ws.Range("G1").Value = "statusId"
$ws.Range("H1").Value = "country"
$ws.Range("I1").Value = "city"

# Update row 2 existing fields
$ws.Range("B2").Value = "Fake Job"
$ws.Range("C2").Value = "It isnt real, you are paid `$150,000 CAD to do nothing."
$ws.Range("D2").Value = "Impossible Inc"
$ws.Range("E2").Value = "2021-10-18T02:55:38.252Z"

# New row 2 fields
$ws.Range("F2").Value = "Sent"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "Canada"
$ws.Range("I2").Value = "Niagara on the Lake"

# New row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Also fake"
$ws.Range("C3").Value = "This is is also not real but it pays nothing "
$ws.Range("D3").Value = "Moo Moo Enterprises"
$ws.Range("E3").Value = "2021-10-18T02:57:10.079Z"
$ws.Range("F3").Value = "Sent"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = "Antarctica"
$ws.Range("I3").Value = "Godrich"
